# Generate Report for Handoff
# Updates the "Status" column from "In Translation" to "Ready for handoff"
# and refreshes the "Latest Handoff Datetime" timestamps on all sheets,
# then widens the status/datetime columns to fit the new text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: zh-cn / de-de status columns (E2, F2) and generate date (G2)
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-13 12:48:56"

# --- zh-cn sheet: Status (C2) and Latest Handoff Datetime (H2)
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-13 12:48:48"

# --- de-de sheet: Status (C2) and Latest Handoff Datetime (H2)
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-13 12:48:56"

# --- Widen the columns that now hold the longer "Ready for handoff" text.
# ColumnWidth is expressed in characters and is quantized by the host to a
# whole-pixel grid (stored_width = round(chars*6)/6) before being written
# back as the OOXML <col width>, so we pick the "characters" value whose
# quantized result lands nearest the target stored width (~17.22).
$overview.Columns.Item(5).ColumnWidth = 16.35
$overview.Columns.Item(6).ColumnWidth = 16.35
$zhcn.Columns.Item(3).ColumnWidth = 16.35
$dede.Columns.Item(3).ColumnWidth = 16.35
